$d = $word.ActiveDocument

$null = $d.Content.Find.Execute("270÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "405÷5=", 2)
$null = $d.Content.Find.Execute("275÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "730÷4=", 2)
$null = $d.Content.Find.Execute("539÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "863÷3=", 2)
$null = $d.Content.Find.Execute("460÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "632÷4=", 2)
$null = $d.Content.Find.Execute("103÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "633÷2=", 2)
$null = $d.Content.Find.Execute("307÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "239÷4=", 2)
$null = $d.Content.Find.Execute("905÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "258÷2=", 2)
$null = $d.Content.Find.Execute("906÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "627÷7=", 2)
$null = $d.Content.Find.Execute("488÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "841÷4=", 2)
$null = $d.Content.Find.Execute("787÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "927÷3=", 2)
$null = $d.Content.Find.Execute("205÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "165÷9=", 2)
$null = $d.Content.Find.Execute("732÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "146÷6=", 2)
$null = $d.Content.Find.Execute("149÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "488÷5=", 2)
$null = $d.Content.Find.Execute("966÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "997÷5=", 2)
$null = $d.Content.Find.Execute("214÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "794÷2=", 2)
$null = $d.Content.Find.Execute("230÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "744÷6=", 2)
$null = $d.Content.Find.Execute("755÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "123÷6=", 2)
$null = $d.Content.Find.Execute("395÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "997÷8=", 2)
$null = $d.Content.Find.Execute("299÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "912÷6=", 2)
$null = $d.Content.Find.Execute("920÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "146÷9=", 2)
$null = $d.Content.Find.Execute("122÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "651÷8=", 2)
$null = $d.Content.Find.Execute("426÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "883÷7=", 2)
$null = $d.Content.Find.Execute("173÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "957÷4=", 2)
$null = $d.Content.Find.Execute("896÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "618÷9=", 2)
$null = $d.Content.Find.Execute("810÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "582÷2=", 2)
